$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 168 (which holds ID 1099001 /
# "Dream Interpretation"). This pushes that existing row down to row 171 and
# leaves rows 168-170 empty, ready to be populated with the three new
# dispatch entries (1052001, 1052002, 1052003).
$ws.Range("A168:A170").EntireRow.Insert()

# The reward-quantity columns (AB, AH) store a quantity that looks numeric
# ("1.0") but must be kept as literal text, matching the source data. A
# leading apostrophe forces Excel to treat the entry as text instead of
# silently coercing it to the number 1.
$qty = "'1.0"

# ---------------------------------------------------------------------
# Row 168 - 1052001 "Unconventional Wedding"
# ---------------------------------------------------------------------
$ws.Range("A168").Value = 1052001
$ws.Range("B168").Value = "Purple"
$ws.Range("C168").Value = "Unconventional Wedding"
$ws.Range("D168").Value = "ニッチな結婚式"
$ws.Range("E168").Value = "독특한 결혼식"
$ws.Range("F168").Value = "小众婚礼"
$ws.Range("G168").Value = "小眾婚禮"
$ws.Range("H168").Value = "A rich heiress from Eastside known for her unconventional hobbies has recently gotten engaged and is planning a paranormal-themed wedding. This has caused many wedding planners to back away from the project. When the news reaches the Bureau, a certain Sinner seems very intrigued."
$ws.Range("I168").Value = "ニューシティのとある名家の令嬢は、ニッチな趣味を持つことで知られている。彼女は最近婚約したばかりで、幽霊をテーマにした結婚式を希望しており、多くのウェディングプランナーが尻込みしている。その噂が管理局に届くと、あるコンビクトが強い関心を示した。"
$ws.Range("J168").Value = "신성의 어느 귀족 집안 아가씨는 독특한 취미로 이름이 높았다. 최근 약혼한 그녀는 오컬트 테마의 결혼식을 원했고, 이는 수많은 웨딩 플래너들을 난처하게 했다. 이 소식을 들은 한 수감자가 큰 흥미를 보였다."
$ws.Range("K168").Value = "新城某贵族千金以爱好小众而闻名，她刚刚订婚，计划为自己办一场灵异主题的婚礼，这叫许多婚礼策划人望而却步。消息传到管理局，某位禁闭者表现出了不小的兴趣。"
$ws.Range("L168").Value = "新城某貴族千金以愛好小眾而聞名，她剛剛訂婚，計畫為自己辦一場靈異主題的婚禮，這叫許多婚禮企劃人望而卻步。消息傳到管理局，某位禁閉者表現出了不小的興趣。"
$ws.Range("M168").Value = "Graves"
$ws.Range("N168").Value = "グレイヴ"
$ws.Range("O168").Value = "그레이브"
$ws.Range("P168").Value = "格芮芙"
$ws.Range("Q168").Value = "格芮芙"
$ws.Range("W168").Value = "Infected Elytra"
$ws.Range("X168").Value = "感染鞘翅"
$ws.Range("Y168").Value = "오염된 겉날개"
$ws.Range("Z168").Value = "感染鞘翅"
$ws.Range("AA168").Value = "感染鞘翅"
$ws.Range("AB168").Value = $qty
$ws.Range("AC168").Value = "Organic Elytra"
$ws.Range("AD168").Value = "原生鞘翅"
$ws.Range("AE168").Value = "원시적 겉날개"
$ws.Range("AF168").Value = "原生鞘翅"
$ws.Range("AG168").Value = "原生鞘翅"
$ws.Range("AH168").Value = $qty

# ---------------------------------------------------------------------
# Row 169 - 1052002 "Home Makeover"
# ---------------------------------------------------------------------
$ws.Range("A169").Value = 1052002
$ws.Range("B169").Value = "Green"
$ws.Range("C169").Value = "Home Makeover"
$ws.Range("D169").Value = "住宅改造"
$ws.Range("E169").Value = "주택 개조"
$ws.Range("F169").Value = "房屋改造"
$ws.Range("G169").Value = "房屋改造"
$ws.Range("H169").Value = "A real estate company has several haunted houses they simply can't sell off, no matter what they try. They're willing to pay a high price for someone to help resolve this desperate situation. Upon hearing this, a certain Sinner eagerly volunteers."
$ws.Range("I169").Value = "ある不動産会社が、なかなか売れない事故物件を複数抱えており、高額な報酬で対応できる人材を探している。それを聞いたあるコンビクトが積極的に参加を申し出た。"
$ws.Range("J169").Value = "한 부동산 업체가 여러 개의 흉가 매물을 가지고 있었는데, 아무리 노력해도 팔리지 않아 흉가 매물을 처리해 주는 사람에게 후한 보상을 주겠다고 약속했다. 이 소식을 들은 한 수감자가 적극적으로 이에 지원했다."
$ws.Range("K169").Value = "某房地产公司手头有好几间凶宅，无论如何也抛售不出，希望有人能为他们解决燃眉之急，奖金丰厚。某禁闭者听说后积极要求参与。"
$ws.Range("L169").Value = "某房地產公司手頭有好幾間凶宅，無論如何也拋售不出，希望有人能為他們解決燃眉之急，獎金豐厚。某禁閉者聽說後積極要求參與。"
$ws.Range("M169").Value = "Luminita"
$ws.Range("N169").Value = "ルミニタ"
$ws.Range("O169").Value = "루미니타"
$ws.Range("P169").Value = "卢米尼塔"
$ws.Range("Q169").Value = "盧米尼塔"
$ws.Range("W169").Value = "Organic Elytra"
$ws.Range("X169").Value = "原生鞘翅"
$ws.Range("Y169").Value = "원시적 겉날개"
$ws.Range("Z169").Value = "原生鞘翅"
$ws.Range("AA169").Value = "原生鞘翅"
$ws.Range("AB169").Value = $qty
$ws.Range("AC169").Value = "Elytra Shard"
$ws.Range("AD169").Value = "鞘翅の破片"
$ws.Range("AE169").Value = "겉날개 파편"
$ws.Range("AF169").Value = "鞘翅残片"
$ws.Range("AG169").Value = "鞘翅殘片"
$ws.Range("AH169").Value = $qty

# ---------------------------------------------------------------------
# Row 170 - 1052003 "Dis Haunt"
# ---------------------------------------------------------------------
$ws.Range("A170").Value = 1052003
$ws.Range("B170").Value = "Blue"
$ws.Range("C170").Value = "Dis Haunt"
$ws.Range("D170").Value = "戦慄の映画製作"
$ws.Range("E170").Value = "공포 전문가"
$ws.Range("F170").Value = "月光光心慌慌"
$ws.Range("G170").Value = "月光光心慌慌"
$ws.Range("H170").Value = "A prominent Eastside director has just announced a new horror film project. They are now urgently recruiting consultants savvy in paranormal phenomena to join the crew."
$ws.Range("I170").Value = "最近、ニューシティのとある有名監督の最新ホラー映画が準備段階の初期に入ったと発表された。現在、撮影協力として霊的分野の専門アドバイザーを多数急募している。"
$ws.Range("J170").Value = "최근 신성의 한 유명 감독이 차기 공포 영화 제작을 공식 발표했고 현재 기괴한 현상에 대해 잘 알고 있는 전문가를 급히 모집하고 있다."
$ws.Range("K170").Value = "近日，新城某知名导演的最新恐怖大片宣布进入早期筹备阶段，现正急招多名灵异方面的专业顾问协助拍摄。"
$ws.Range("L170").Value = "近日，新城某知名導演的最新恐怖大片宣佈進入早期籌備階段，現正急招多名靈異方面的專業顧問協助拍攝。"
$ws.Range("M170").Value = "Luminita"
$ws.Range("N170").Value = "ルミニタ"
$ws.Range("O170").Value = "루미니타"
$ws.Range("P170").Value = "卢米尼塔"
$ws.Range("Q170").Value = "盧米尼塔"
$ws.Range("R170").Value = "Letta"
$ws.Range("S170").Value = "レタ"
$ws.Range("T170").Value = "레타"
$ws.Range("U170").Value = "莱塔"
$ws.Range("V170").Value = "萊塔"
$ws.Range("W170").Value = "Organic Elytra"
$ws.Range("X170").Value = "原生鞘翅"
$ws.Range("Y170").Value = "원시적 겉날개"
$ws.Range("Z170").Value = "原生鞘翅"
$ws.Range("AA170").Value = "原生鞘翅"
$ws.Range("AB170").Value = $qty
$ws.Range("AC170").Value = "Organic Elytra"
$ws.Range("AD170").Value = "原生鞘翅"
$ws.Range("AE170").Value = "원시적 겉날개"
$ws.Range("AF170").Value = "原生鞘翅"
$ws.Range("AG170").Value = "原生鞘翅"
$ws.Range("AH170").Value = $qty
